$wb = $excel.ActiveWorkbook

# Shared text used on both language sheets.
$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/oltest/blob/981bcaffd22a0a4b75315007f13f7b9d41b3abe0/e2e/48188f1e-5572-431b-b342-dc7310300173.md, latest: https://github.com/OpenLocalizationTestOrg/oltest/blob/5e91bb8510101be45987d25fbff177c38bd14471/e2e/48188f1e-5572-431b-b342-dc7310300173.md."
$currentMdUrl = "https://github.com/OpenLocalizationTestOrg/oltest/blob/981bcaffd22a0a4b75315007f13f7b9d41b3abe0/e2e/48188f1e-5572-431b-b342-dc7310300173.md"
$currentMdDisplay = "48188f1e-5572-431b-b342-dc7310300173.md"

# ---------- zh-cn sheet (row 6: 48188f1e-... handback report) ----------
$wsZh = $wb.Worksheets.Item("zh-cn")

# Latest Target File now links back to the current (out of date) handback source commit.
$wsZh.Hyperlinks.Add($wsZh.Range("I6"), $currentMdUrl, "", "", $currentMdDisplay)

# Latest Handback File / Latest Handback DateTime / Error Detail get populated.
$wsZh.Range("J6").Value = "48188f1e-5572-431b-b342-dc7310300173.41405433557c7e6d931341897ecfbb96a76614ff.zh-cn.xlf"
$wsZh.Range("K6").Value = "2016-08-12 03:04:39"
$wsZh.Range("P6").Value = $errorDetail

# Error Detail column is widened so the long message is readable.
$wsZh.Columns.Item(16).ColumnWidth = 39.17

# ---------- de-de sheet (row 6: 48188f1e-... handback report) ----------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Hyperlinks.Add($wsDe.Range("I6"), $currentMdUrl, "", "", $currentMdDisplay)

$wsDe.Range("J6").Value = "48188f1e-5572-431b-b342-dc7310300173.41405433557c7e6d931341897ecfbb96a76614ff.de-de.xlf"
$wsDe.Range("K6").Value = "2016-08-12 03:04:46"
$wsDe.Range("P6").Value = $errorDetail

$wsDe.Columns.Item(16).ColumnWidth = 39.17
